$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.277.23"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.628.42"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'602.01"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'153.28"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.558"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").Value = "2.626.91"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +5.77%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'27.84"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "3.108.37"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "67.318.96"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "2.588.81"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'363.19"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").Value = "'66.15"
$ws.Range("E26").Value = "  -7.48%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.03"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.759.83"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000104"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "'577.88"
$ws.Range("E30").Value = "  -7.26%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "'157.28"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "'5.27"
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").Value = "'41.21"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D47").Value = "'155.79"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "0.0₆0287"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").Value = "'20.89"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  +0.32%  "

Write-Output "Applied cryptos update"
